$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $oldText"
    }
    $rng.Text = $newText
}

# ---- Title ----
Replace-ExactText "Unveiling the Enigmatic Quantum Realm" "The Fascinating World of Biology: An Exploration of Life's Complexities"

# ---- Author ----
Replace-ExactText "Alex Wissner-Gross" "Emma Johnson"

# ---- Email address (paragraph 3) ----
Replace-ExactText "awissnergross@mit" "emma"
Replace-ExactText "edu" "johnson123@schoolmail"

# Append ".org" as two new runs at the end of the email paragraph (after "schoolmail")
$p3 = $d.Paragraphs(3)
$p3Rng = $p3.Range
$p3Rng.InsertAfter(".")
$p3Rng2 = $p3.Range
$p3Rng2.InsertAfter("org")

# ---- Body paragraph (paragraph 5) ----
Replace-ExactText "In the realm of physics, quantum mechanics stands as a profound and enigmatic chapter, transcending the familiar world of classical physics" "Biology is an intriguing discipline that delves into the intricacies of life, unraveling the secrets of living organisms"

Replace-ExactText " This extraordinary theory delves into the realm of subatomic particles, where the laws of nature are governed by uncertainty and probability" " This vast field encompasses a multitude of interconnected phenomena, ranging from the basic building blocks of life, such as cells, to complex ecological interactions that shape entire ecosystems"

Replace-ExactText " Quantum mechanics has revolutionized our understanding of the universe, unveiling an intricate tapestry of phenomena that defy intuition, from the inexplicable behavior of electrons to the perplexing phenomenon of quantum entanglement" " Biology stands as a testament to the wonders of the natural world, inviting us to explore and understand the remarkable diversity of life on Earth"

Replace-ExactText "As we venture deeper into the quantum realm, we encounter a fascinating paradox" "Biology unveils the intricate machinery of life, from the molecular dances within cells to the intricate web of interactions that govern the functioning of organisms"

Replace-ExactText " The very act of observing quantum systems alters their behavior, rendering them both elusive and unpredictable" " It delves into the study of genetics, the blueprint of life that holds the key to our existence, unraveling the mysteries of inheritance and evolution"

Replace-ExactText " This enigmatic duality, known as the observer effect, has fueled debates among physicists for decades, raising fundamental questions about the nature of reality and the relationship between consciousness and the physical world" " Biology opens up a realm of discovery, where we can marvel at the intricate adaptations of organisms, the delicate balance of ecosystems, and the remarkable resilience of life in the face of adversity"

Replace-ExactText "The profound implications of quantum mechanics extend beyond the theoretical realm, reaching into diverse fields such as computation, cryptography, and medicine" "This captivating discipline invites us to engage in scientific inquiry, encouraging us to pose questions, gather evidence, and analyze data to unravel the mysteries of life"

Replace-ExactText " Quantum computers hold the promise of exponential speed-ups in problem-solving, while quantum cryptography offers unbreakable encryption methods" " It nurtures critical thinking skills, enabling us to evaluate information thoughtfully and make informed decisions"

Replace-ExactText " Additionally, quantum physics is paving the way for advancements in medical imaging, sensing, and drug design" " Biology fosters a deep appreciation for the interconnectedness of life, promoting responsible stewardship of our planet and its inhabitants"

# ---- Summary heading paragraph text (paragraph 7) ----
Replace-ExactText "This essay has provided a glimpse into the enigmatic quantum realm, exploring the profound implications of quantum mechanics on our understanding of the universe" "Biology presents a fascinating journey into the realm of living organisms, unraveling the intricacies of life's complexities"

Replace-ExactText " From the strange world of subatomic particles to its transformative applications across multiple disciplines, quantum mechanics continues to challenge our perceptions of reality and redefine the boundaries of human knowledge" " It encompasses a wide range of topics, including cell biology, genetics, ecology, and evolution"

# Insert three additional sentences (with their own period runs) before the
# paragraph's final, unchanged "." run.
$tail = $d.Content
$tail.Find.Execute("It encompasses a wide range of topics, including cell biology, genetics, ecology, and evolution", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Collapse(0)
$tail.InsertBefore(" This captivating discipline invites us to explore the wonders of the natural world, revealing the remarkable diversity of life on Earth and inspiring us to become responsible stewards of our planet")
$tail.InsertBefore(".")
$tail.InsertBefore(" Biology not only imparts knowledge but also nurtures critical thinking skills, fostering a deep appreciation for the interconnectedness of life")
$tail.InsertBefore(".")

# ---- New trailing empty paragraph at the end of the document body ----
$endOfDoc = $d.Range($d.Content.End, $d.Content.End)
$endOfDoc.Text = "`r"
